$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed column P (rows 1-11) with column O's formatting (number format, borders,
# alignment) by copying O1:O11 -> P1:P11. The actual per-row values below then
# overwrite the copied O-column values, but the cell styles (date format for
# row 1, centered/bordered for the rest) stay intact and reuse the existing
# style entries instead of minting new ones.
$ws.Range("O1:O11").Copy($ws.Range("P1:P11"))

# Also widen column P to match column O (both are week columns).
$ws.Range("P1").EntireColumn.ColumnWidth = $ws.Range("O1").EntireColumn.ColumnWidth

# Week-of 2015-12-14 (15th week) header date.
$ws.Range("P1").Value = 42352

# Per-student scores for the new week.
$ws.Range("P2").Value = 5
$ws.Range("P3").Value = 5
$ws.Range("P4").Value = 5
$ws.Range("P5").Value = 5
$ws.Range("P6").Value = 2
$ws.Range("P7").Value = 5
$ws.Range("P8").Value = 5
$ws.Range("P9").Value = 5
$ws.Range("P10").Value = 5

# Row 11 (totals/legend row) gets "-" in both the existing O column and the
# freshly added P column.
$ws.Range("O11").Value = "-"
$ws.Range("P11").Value = "-"

# Move the active selection to match the author's final cursor position.
$ws.Range("O12").Select()
